# Weekly refresh: a new price entry was recorded for "Poroto granado" at
# Feria Lagunitas de Puerto Montt. It belongs at the top of the date-sorted
# block (row 10), so push the existing rows 10:43 down to 11:44 and fill in
# the new week's data in the freed-up row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 10:43 down to 11:44 (mirrors an Excel "insert row" above row 10).
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with this week's record. Columns A, B, C, E, F, G,
# H, N, Q, R are constant across every row in this sheet.
$ws.Cells.Item(10, 1).Value = 4
$ws.Cells.Item(10, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(10, 3).Value = "Los Lagos"
$ws.Cells.Item(10, 4).Value = "2022-04-08"
$ws.Cells.Item(10, 5).Value = 10
$ws.Cells.Item(10, 6).Value = 100112030
$ws.Cells.Item(10, 7).Value = "Poroto granado"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 60
$ws.Cells.Item(10, 11).Value = 25000
$ws.Cells.Item(10, 12).Value = 25000
$ws.Cells.Item(10, 13).Value = 25000
$ws.Cells.Item(10, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(10, 15).Value = "Región Metropolitana"
$ws.Cells.Item(10, 16).Value = 1000
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
